$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    paragraph (the very first paragraph in the body).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)

# Grab a style-less paragraph's FormattedText (paragraph 3 - a plain body
# paragraph with no pPr/style) and drop a copy of it right after the title.
# This gives us a brand-new paragraph that has no paragraph style applied
# (matching the "Normal"/no-pPr look of the rest of the body) without
# touching .Style (which stamps rsid bookkeeping attributes onto the <w:p>).
$donorPara = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertPoint.FormattedText = $donorPara.Range.FormattedText

$newPara = $d.Paragraphs.Item(2)
$newRange = $newPara.Range
$newRange.End = $newRange.End - 1   # exclude the paragraph mark
$newRange.Text = ""                 # wipe the donor's text, keep the (styleless) paragraph

$boldPart = "Meta description"
$restPart = ": Explore the Wild West and find out the pros and cons of Cowboy Treasure Deluxe online slot. Play now for free."
$fullText = $boldPart + $restPart

$typePoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$typePoint.InsertAfter($fullText)

$boldRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $boldPart.Length)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph
#    ("Play Cowboy Treasure Deluxe Free - Pros and Cons | Review")
#    that sits near the end of the document, right before the italic
#    meta-description paragraph. (Skip paragraph 1, which is the real
#    Heading-1 document title and must stay.)
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 2; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Cowboy Treasure Deluxe Free - Pros and Cons | Review`r") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 3) Replace the old meta-description sentence (now the very last paragraph,
#    in italics) with the new AI image-generation prompt text. Scope the
#    Find to the last paragraph's own Range so it can't accidentally match
#    the "Meta description: Explore the Wild West..." text we just added
#    near the top of the document.
# ---------------------------------------------------------------------------
$newPrompt = 'Prompt: Create a cartoon-style feature image for Cowboy Treasure Deluxe that features a happy Maya warrior with glasses. For the feature image, we want to blend the Wild West theme with a fun twist. We want to feature a cartoon-style image of a happy Maya warrior with glasses, holding up a gold nugget with one hand and a sheriff star with the other hand. The Maya warrior should have a big smile on his face and be wearing traditional Maya clothing with a cowboy hat. The background should feature the reddish rock formations that resemble the Grand Canyon, just like in the game. The image should also include the game title "Cowboy Treasure Deluxe" written in a fun Western-style font. Overall, the image should convey the excitement and joy of playing Cowboy Treasure Deluxe with a touch of humor.'

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.Find.Execute("Explore the Wild West and find out the pros and cons of Cowboy Treasure Deluxe online slot. Play now for free.", $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2)

Write-Host "Paragraph count: $($d.Paragraphs.Count)"
